# Update the "nao_respondidos_para_ligar" sheet:
# - Remove the old data rows (2-6), keeping only the header row
# - Change the header row from NOME/TELEFONE to NOME/NUMERO/valor/vencimento/status/STATUS

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old sample data rows (rows 2 through 6)
$ws.Range("A2:B6").EntireRow.Delete()

# Update/extend the header row
$ws.Range("A1").Value = "NOME"
$ws.Range("B1").Value = "NUMERO"
$ws.Range("C1").Value = "valor"
$ws.Range("D1").Value = "vencimento"
$ws.Range("E1").Value = "status"
$ws.Range("F1").Value = "STATUS"

# Apply the existing header style (bold, bordered, centered) to the new header cells
$ws.Range("B1").Copy()
$ws.Range("C1:F1").PasteSpecial(-4122)
